$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 731.3
$ws.Range("I33").Value = 259
$ws.Range("K33").Value = 259
$ws.Range("M33").Value = -30

# Row 51
$ws.Range("H51").Value = 2978.5715
$ws.Range("J51").Value = 2978.5715
$ws.Range("L51").Value = 2978.5715
$ws.Range("N51").Value = -3946.5715

# Row 69
$ws.Range("H69").Value = 6965.25
$ws.Range("J69").Value = 7096.2666
$ws.Range("L69").Value = 21288.7998
$ws.Range("N69").Value = -23036.7998

# Row 72
$ws.Range("H72").Value = 6965.25
$ws.Range("J72").Value = 7096.2666
$ws.Range("L72").Value = 63866.3994
$ws.Range("N72").Value = -72602.39939999999

# Row 86
$ws.Range("H86").Value = 4596.2
$ws.Range("I86").Value = 4596.2
$ws.Range("K86").Value = 4596.2
$ws.Range("M86").Value = -3473.2

# Row 89
$ws.Range("H89").Value = 4596.2
$ws.Range("I89").Value = 4596.2
$ws.Range("K89").Value = 22981
$ws.Range("M89").Value = -17365

# Row 113
$ws.Range("H113").Value = 5618.75
$ws.Range("I113").Value = 2987.5
$ws.Range("K113").Value = 2987.5
$ws.Range("M113").Value = 266.5

# Row 127
$ws.Range("H127").Value = 2183.5715
$ws.Range("I127").Value = 1047.5
$ws.Range("K127").Value = 3142.5
$ws.Range("M127").Value = 1817.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2144.4375
$ws.Range("I45").Value = 1763.3077
$ws.Range("K45").Value = 1763.3077
$ws.Range("M45").Value = -1386.3077

# Row 61
$ws.Range("H61").Value = 4088.9
$ws.Range("I61").Value = 2398.25
$ws.Range("J61").Value = 6624.875
$ws.Range("K61").Value = 2398.25
$ws.Range("L61").Value = 6624.875
$ws.Range("M61").Value = -2186.25
$ws.Range("N61").Value = -7048.875

# Row 110
$ws.Range("H110").Value = 791
$ws.Range("I110").Value = 810
$ws.Range("J110").Value = 724.5
$ws.Range("K110").Value = 810
$ws.Range("L110").Value = 724.5
$ws.Range("M110").Value = 1235
$ws.Range("N110").Value = -4814.5

# Row 132
$ws.Range("H132").Value = 3099.1538
$ws.Range("I132").Value = 2979.8096
$ws.Range("J132").Value = 3600.4
$ws.Range("K132").Value = 8939.4288
$ws.Range("L132").Value = 10801.2
$ws.Range("M132").Value = -6409.4288
$ws.Range("N132").Value = -15861.2

# Row 136
$ws.Range("H136").Value = 4088.9
$ws.Range("I136").Value = 2398.25
$ws.Range("J136").Value = 6624.875
$ws.Range("K136").Value = 7194.75
$ws.Range("L136").Value = 19874.625
$ws.Range("M136").Value = -4644.75
$ws.Range("N136").Value = -24974.625

$ws = $wb.Worksheets.Item("BSM")
# Row 87
$ws.Range("H87").Value = 90000
$ws.Range("J87").Value = 90000
$ws.Range("L87").Value = 90000
$ws.Range("N87").Value = -92496

# Row 90
$ws.Range("H90").Value = 90000
$ws.Range("J90").Value = 90000
$ws.Range("L90").Value = 270000
$ws.Range("N90").Value = -282480

# Row 94
$ws.Range("H94").Value = 1794.8462
$ws.Range("I94").Value = 1812.1818
$ws.Range("K94").Value = 1812.1818
$ws.Range("M94").Value = -1361.1818

# Row 105
$ws.Range("H105").Value = 3079.7334
$ws.Range("I105").Value = 3056.6365
$ws.Range("J105").Value = 3143.25
$ws.Range("K105").Value = 3056.6365
$ws.Range("L105").Value = 3143.25
$ws.Range("M105").Value = -1309.6365
$ws.Range("N105").Value = -6637.25

# Row 134
$ws.Range("H134").Value = 1549.125
$ws.Range("I134").Value = 1519.0667
$ws.Range("K134").Value = 4557.2001
$ws.Range("M134").Value = -2022.2001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6074.4136
$ws.Range("I31").Value = 1636.5
$ws.Range("J31").Value = 7765.048
$ws.Range("K31").Value = 1636.5
$ws.Range("L31").Value = 7765.048
$ws.Range("M31").Value = -1341.5
$ws.Range("N31").Value = -8355.047999999999

# Row 34
$ws.Range("H34").Value = 6074.4136
$ws.Range("I34").Value = 1636.5
$ws.Range("J34").Value = 7765.048
$ws.Range("K34").Value = 1636.5
$ws.Range("L34").Value = 7765.048
$ws.Range("M34").Value = -1434.5
$ws.Range("N34").Value = -8169.048

# Row 56
$ws.Range("H56").Value = 4000
$ws.Range("I56").Value = 4000
$ws.Range("K56").Value = 4000
$ws.Range("M56").Value = -3155

# Row 86
$ws.Range("H86").Value = 6902.3335
$ws.Range("I86").Value = 5349.5
$ws.Range("K86").Value = 5349.5
$ws.Range("M86").Value = -4226.5

# Row 89
$ws.Range("H89").Value = 6902.3335
$ws.Range("I89").Value = 5349.5
$ws.Range("K89").Value = 26747.5
$ws.Range("M89").Value = -21131.5

# Row 99
$ws.Range("H99").Value = 1885.4
$ws.Range("I99").Value = 1575.6666
$ws.Range("J99").Value = 2350
$ws.Range("K99").Value = 1575.6666
$ws.Range("L99").Value = 2350
$ws.Range("M99").Value = -77.66660000000002
$ws.Range("N99").Value = -5346

# Row 122
$ws.Range("H122").Value = 2333.3333
$ws.Range("I122").Value = 2125
$ws.Range("K122").Value = 6375
$ws.Range("M122").Value = -3925

# Row 126
$ws.Range("H126").Value = 1885.4
$ws.Range("I126").Value = 1575.6666
$ws.Range("J126").Value = 2350
$ws.Range("K126").Value = 4726.9998
$ws.Range("L126").Value = 7050
$ws.Range("M126").Value = -2256.9998
$ws.Range("N126").Value = -11990

# Row 134
$ws.Range("H134").Value = 1111.75
$ws.Range("J134").Value = 600
$ws.Range("L134").Value = 1800
$ws.Range("N134").Value = -6870

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 159835.38
$ws.Range("J37").Value = 159835.38
$ws.Range("L37").Value = 479506.14
$ws.Range("N37").Value = -479730.14

# Row 58
$ws.Range("H58").Value = 1133
$ws.Range("J58").Value = 699.5
$ws.Range("L58").Value = 2098.5
$ws.Range("N58").Value = -2354.5

# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null

$ws = $wb.Worksheets.Item("GSM")
# Row 6
$ws.Range("H6").Value = 5000
$ws.Range("J6").Value = 5000
$ws.Range("L6").Value = 5000
$ws.Range("N6").Value = -5226

# Row 14
$ws.Range("H14").Value = 3560.6428
$ws.Range("J14").Value = 6361.25
$ws.Range("L14").Value = 6361.25
$ws.Range("N14").Value = -6697.25

# Row 16
$ws.Range("H16").Value = 5000
$ws.Range("J16").Value = 5000
$ws.Range("L16").Value = 5000
$ws.Range("N16").Value = -5500

# Row 27
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = $null

# Row 46
$ws.Range("H46").Value = 14500
$ws.Range("J46").Value = 15000
$ws.Range("L46").Value = 15000
$ws.Range("N46").Value = -15312

# Row 126
$ws.Range("H126").Value = 3998.75
$ws.Range("I126").Value = 2995
$ws.Range("K126").Value = 8985
$ws.Range("M126").Value = -6515

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6860.846
$ws.Range("I7").Value = 6239
$ws.Range("K7").Value = 6239
$ws.Range("M7").Value = -6127

# Row 46
$ws.Range("H46").Value = 3108.6667
$ws.Range("I46").Value = 3280.4
$ws.Range("J46").Value = 2250
$ws.Range("K46").Value = 3280.4
$ws.Range("L46").Value = 2250
$ws.Range("M46").Value = -3092.4
$ws.Range("N46").Value = -2626

# Row 82
$ws.Range("H82").Value = 8000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 8000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 8000
$ws.Range("M82").Value = $null
$ws.Range("N82").Value = -8722

# Row 85
$ws.Range("H85").Value = 8000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 8000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 8000
$ws.Range("M85").Value = $null
$ws.Range("N85").Value = -10496

# Row 126
$ws.Range("H126").Value = 6860.846
$ws.Range("I126").Value = 6239
$ws.Range("K126").Value = 18717
$ws.Range("M126").Value = -16247

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2562.4443
$ws.Range("I132").Value = 2533.6667
$ws.Range("J132").Value = 2620
$ws.Range("K132").Value = 7601.000100000001
$ws.Range("L132").Value = 7860
$ws.Range("M132").Value = -5071.000100000001
$ws.Range("N132").Value = -12920
